$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark near the start of the document
#    (it previously sat between the "2" run and the following space run).
$d.Bookmarks("_GoBack").Delete()

# 2. Split the paragraph that holds the "_Hlk536537895" bookmark right after
#    "...to be pre-installed." so that bookmark end stays with the first
#    half and a new paragraph begins right after it.
$bm = $d.Bookmarks("_Hlk536537895")
$bmRange = $bm.Range
$splitPoint = $bmRange.End

$insertionRange = $d.Range($splitPoint, $splitPoint)
$insertionRange.InsertParagraphAfter()

# Re-anchor the "_Hlk536537895" bookmark so its end mark stays inside the
# now-closed first paragraph (InsertParagraphAfter pushes bookmarks that sit
# exactly at the split point into the new paragraph otherwise).
$restoredRange = $d.Range($bmRange.Start, $splitPoint)
$d.Bookmarks.Add("_Hlk536537895", $restoredRange)

# 3. Re-create the "_GoBack" bookmark at the very start of the new paragraph.
$goBackRange = $d.Range($splitPoint + 1, $splitPoint + 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
